# Prefix every colon-separated hex byte in the "doip" (G) and "uds" (H)
# columns with "0x" so downstream lookups can match on the 0x-prefixed
# codes (e.g. "02:fd:00" -> "0x02:0xfd:0x00"). Cells holding the literal
# "N/A" placeholder are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 41

for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne "N/A" -and $val -ne "") {
            $parts = $val -split ":"
            $newVal = ($parts | ForEach-Object { "0x" + $_ }) -join ":"
            $cell.Value = $newVal
        }
    }
}
